$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("D1").Value = "CRIT_RATE"
$ws.Range("E1").Value = "CRIT_DAMAGE"

# Data rows
$critRate = @(5, 10, 10, 15, 15, 20, 25, 30, 30)
$critDamage = @(10, 10, 20, 20, 30, 30, 40, 40, 50)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $critRate[$i]
    $ws.Cells.Item($row, 5).Value = $critDamage[$i]
}

$ws.Range("F4").Select()
